$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date/time number format + style from an existing column-A data cell
# so the new rows reuse the same cellXf (s="2") instead of minting a new style.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A1169:A1232").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$data = New-Object "object[,]" 64,6
$data[0,0] = 45534.5
$data[0,1] = 0.3513
$data[0,2] = 0.3546
$data[0,3] = 0.3461
$data[0,4] = 0.3519
$data[0,5] = 24333063.7
$data[1,0] = 45534.58333333334
$data[1,1] = 0.352
$data[1,2] = 0.3543
$data[1,3] = 0.3383
$data[1,4] = 0.34
$data[1,5] = 23508445.5
$data[2,0] = 45534.66666666666
$data[2,1] = 0.34
$data[2,2] = 0.3436
$data[2,3] = 0.3373
$data[2,4] = 0.343
$data[2,5] = 10650201.5
$data[3,0] = 45534.75
$data[3,1] = 0.3429
$data[3,2] = 0.3493
$data[3,3] = 0.3426
$data[3,4] = 0.3455
$data[3,5] = 5988210.1
$data[4,0] = 45534.83333333334
$data[4,1] = 0.3455
$data[4,2] = 0.3479
$data[4,3] = 0.3455
$data[4,4] = 0.3469
$data[4,5] = 1846793.3
$data[5,0] = 45534.91666666666
$data[5,1] = 0.3469
$data[5,2] = 0.3482
$data[5,3] = 0.3465
$data[5,4] = 0.3469
$data[5,5] = 2579506.8
$data[6,0] = 45535
$data[6,1] = 0.3469
$data[6,2] = 0.3485
$data[6,3] = 0.3458
$data[6,4] = 0.3466
$data[6,5] = 2661263.8
$data[7,0] = 45535.08333333334
$data[7,1] = 0.3467
$data[7,2] = 0.3477
$data[7,3] = 0.3454
$data[7,4] = 0.3475
$data[7,5] = 2132374.7
$data[8,0] = 45535.16666666666
$data[8,1] = 0.3474
$data[8,2] = 0.3504
$data[8,3] = 0.3474
$data[8,4] = 0.348
$data[8,5] = 2663054.2
$data[9,0] = 45535.25
$data[9,1] = 0.348
$data[9,2] = 0.3501
$data[9,3] = 0.3475
$data[9,4] = 0.3478
$data[9,5] = 2498661.6
$data[10,0] = 45535.33333333334
$data[10,1] = 0.3478
$data[10,2] = 0.3484
$data[10,3] = 0.3453
$data[10,4] = 0.3458
$data[10,5] = 2402320.6
$data[11,0] = 45535.41666666666
$data[11,1] = 0.3458
$data[11,2] = 0.3495
$data[11,3] = 0.3454
$data[11,4] = 0.3481
$data[11,5] = 2407327.8
$data[12,0] = 45535.5
$data[12,1] = 0.3481
$data[12,2] = 0.349
$data[12,3] = 0.3473
$data[12,4] = 0.3475
$data[12,5] = 2256718.4
$data[13,0] = 45535.58333333334
$data[13,1] = 0.3475
$data[13,2] = 0.3496
$data[13,3] = 0.3471
$data[13,4] = 0.3492
$data[13,5] = 2455412.3
$data[14,0] = 45535.66666666666
$data[14,1] = 0.3492
$data[14,2] = 0.3493
$data[14,3] = 0.345
$data[14,4] = 0.3471
$data[14,5] = 2296348.7
$data[15,0] = 45535.75
$data[15,1] = 0.3469
$data[15,2] = 0.3486
$data[15,3] = 0.3455
$data[15,4] = 0.3468
$data[15,5] = 2455239.7
$data[16,0] = 45535.83333333334
$data[16,1] = 0.3467
$data[16,2] = 0.347
$data[16,3] = 0.3447
$data[16,4] = 0.3449
$data[16,5] = 1160812.6
$data[17,0] = 45535.91666666666
$data[17,1] = 0.3462
$data[17,2] = 0.3464
$data[17,3] = 0.3444
$data[17,4] = 0.345
$data[17,5] = 1188086.3
$data[18,0] = 45536
$data[18,1] = 0.3451
$data[18,2] = 0.3458
$data[18,3] = 0.3413
$data[18,4] = 0.3414
$data[18,5] = 2784322.5
$data[19,0] = 45536.08333333334
$data[19,1] = 0.3414
$data[19,2] = 0.3449
$data[19,3] = 0.3394
$data[19,4] = 0.3448
$data[19,5] = 3618684.9
$data[20,0] = 45536.16666666666
$data[20,1] = 0.3449
$data[20,2] = 0.3449
$data[20,3] = 0.3393
$data[20,4] = 0.3402
$data[20,5] = 3761037.6
$data[21,0] = 45536.25
$data[21,1] = 0.3401
$data[21,2] = 0.3446
$data[21,3] = 0.34
$data[21,4] = 0.3446
$data[21,5] = 2217267.8
$data[22,0] = 45536.33333333334
$data[22,1] = 0.3446
$data[22,2] = 0.345
$data[22,3] = 0.3429
$data[22,4] = 0.3439
$data[22,5] = 1562092.2
$data[23,0] = 45536.41666666666
$data[23,1] = 0.3439
$data[23,2] = 0.344
$data[23,3] = 0.3403
$data[23,4] = 0.3432
$data[23,5] = 3077041.3
$data[24,0] = 45536.5
$data[24,1] = 0.3431
$data[24,2] = 0.3436
$data[24,3] = 0.3363
$data[24,4] = 0.3405
$data[24,5] = 6828323.3
$data[25,0] = 45536.58333333334
$data[25,1] = 0.3406
$data[25,2] = 0.3416
$data[25,3] = 0.3351
$data[25,4] = 0.3403
$data[25,5] = 8932560.1
$data[26,0] = 45536.66666666666
$data[26,1] = 0.3403
$data[26,2] = 0.3417
$data[26,3] = 0.3341
$data[26,4] = 0.3348
$data[26,5] = 7220129.1
$data[27,0] = 45536.75
$data[27,1] = 0.3349
$data[27,2] = 0.3422
$data[27,3] = 0.3347
$data[27,4] = 0.3397
$data[27,5] = 6474757.6
$data[28,0] = 45536.83333333334
$data[28,1] = 0.3397
$data[28,2] = 0.3421
$data[28,3] = 0.3357
$data[28,4] = 0.3373
$data[28,5] = 3231814.7
$data[29,0] = 45536.91666666666
$data[29,1] = 0.3374
$data[29,2] = 0.3378
$data[29,3] = 0.3267
$data[29,4] = 0.3311
$data[29,5] = 14158076.3
$data[30,0] = 45537
$data[30,1] = 0.331
$data[30,2] = 0.3345
$data[30,3] = 0.3308
$data[30,4] = 0.3331
$data[30,5] = 3754763.2
$data[31,0] = 45537.08333333334
$data[31,1] = 0.3332
$data[31,2] = 0.3345
$data[31,3] = 0.331
$data[31,4] = 0.3316
$data[31,5] = 2886182.4
$data[32,0] = 45537.16666666666
$data[32,1] = 0.3315
$data[32,2] = 0.3316
$data[32,3] = 0.3271
$data[32,4] = 0.3301
$data[32,5] = 5411423.9
$data[33,0] = 45537.25
$data[33,1] = 0.3301
$data[33,2] = 0.3312
$data[33,3] = 0.325
$data[33,4] = 0.3259
$data[33,5] = 5903645.7
$data[34,0] = 45537.33333333334
$data[34,1] = 0.3259
$data[34,2] = 0.3327
$data[34,3] = 0.3253
$data[34,4] = 0.3322
$data[34,5] = 13824019.5
$data[35,0] = 45537.41666666666
$data[35,1] = 0.3322
$data[35,2] = 0.3364
$data[35,3] = 0.3314
$data[35,4] = 0.3346
$data[35,5] = 5705661.6
$data[36,0] = 45537.5
$data[36,1] = 0.3338
$data[36,2] = 0.3349
$data[36,3] = 0.3287
$data[36,4] = 0.3287
$data[36,5] = 5483505.8
$data[37,0] = 45537.58333333334
$data[37,1] = 0.3287
$data[37,2] = 0.3339
$data[37,3] = 0.3269
$data[37,4] = 0.3325
$data[37,5] = 9903972.6
$data[38,0] = 45537.66666666666
$data[38,1] = 0.3325
$data[38,2] = 0.3342
$data[38,3] = 0.3297
$data[38,4] = 0.3332
$data[38,5] = 3418746.4
$data[39,0] = 45537.75
$data[39,1] = 0.3333
$data[39,2] = 0.3335
$data[39,3] = 0.3298
$data[39,4] = 0.3309
$data[39,5] = 5888111.4
$data[40,0] = 45537.83333333334
$data[40,1] = 0.3309
$data[40,2] = 0.3374
$data[40,3] = 0.3303
$data[40,4] = 0.3364
$data[40,5] = 4437607.1
$data[41,0] = 45537.91666666666
$data[41,1] = 0.3364
$data[41,2] = 0.337
$data[41,3] = 0.3347
$data[41,4] = 0.3362
$data[41,5] = 3572685
$data[42,0] = 45538
$data[42,1] = 0.3362
$data[42,2] = 0.337
$data[42,3] = 0.3342
$data[42,4] = 0.3361
$data[42,5] = 4897524.4
$data[43,0] = 45538.08333333334
$data[43,1] = 0.336
$data[43,2] = 0.3398
$data[43,3] = 0.3348
$data[43,4] = 0.3351
$data[43,5] = 5103906.9
$data[44,0] = 45538.16666666666
$data[44,1] = 0.3352
$data[44,2] = 0.3356
$data[44,3] = 0.3329
$data[44,4] = 0.3341
$data[44,5] = 3223761.9
$data[45,0] = 45538.25
$data[45,1] = 0.3341
$data[45,2] = 0.3341
$data[45,3] = 0.3318
$data[45,4] = 0.3334
$data[45,5] = 3787614.8
$data[46,0] = 45538.33333333334
$data[46,1] = 0.3335
$data[46,2] = 0.3336
$data[46,3] = 0.3281
$data[46,4] = 0.3291
$data[46,5] = 6665063.5
$data[47,0] = 45538.41666666666
$data[47,1] = 0.329
$data[47,2] = 0.3312
$data[47,3] = 0.3286
$data[47,4] = 0.3302
$data[47,5] = 3695518
$data[48,0] = 45538.5
$data[48,1] = 0.3302
$data[48,2] = 0.332
$data[48,3] = 0.3227
$data[48,4] = 0.3243
$data[48,5] = 12347933.7
$data[49,0] = 45538.58333333334
$data[49,1] = 0.3242
$data[49,2] = 0.3242
$data[49,3] = 0.32
$data[49,4] = 0.3204
$data[49,5] = 14600325.5
$data[50,0] = 45538.66666666666
$data[50,1] = 0.3204
$data[50,2] = 0.3253
$data[50,3] = 0.3195
$data[50,4] = 0.3233
$data[50,5] = 5665410.9
$data[51,0] = 45538.75
$data[51,1] = 0.3233
$data[51,2] = 0.3259
$data[51,3] = 0.3221
$data[51,4] = 0.3249
$data[51,5] = 3813726
$data[52,0] = 45538.83333333334
$data[52,1] = 0.3249
$data[52,2] = 0.3262
$data[52,3] = 0.3222
$data[52,4] = 0.3231
$data[52,5] = 2980531.7
$data[53,0] = 45538.91666666666
$data[53,1] = 0.3231
$data[53,2] = 0.3236
$data[53,3] = 0.3172
$data[53,4] = 0.3185
$data[53,5] = 8693428.1
$data[54,0] = 45539
$data[54,1] = 0.3184
$data[54,2] = 0.3202
$data[54,3] = 0.305
$data[54,4] = 0.3143
$data[54,5] = 34461923.4
$data[55,0] = 45539.08333333334
$data[55,1] = 0.3143
$data[55,2] = 0.317
$data[55,3] = 0.3135
$data[55,4] = 0.3157
$data[55,5] = 10320149.7
$data[56,0] = 45539.16666666666
$data[56,1] = 0.3157
$data[56,2] = 0.3172
$data[56,3] = 0.3124
$data[56,4] = 0.313
$data[56,5] = 6063714
$data[57,0] = 45539.25
$data[57,1] = 0.313
$data[57,2] = 0.3222
$data[57,3] = 0.313
$data[57,4] = 0.3212
$data[57,5] = 10501210.7
$data[58,0] = 45539.33333333334
$data[58,1] = 0.3212
$data[58,2] = 0.325
$data[58,3] = 0.321
$data[58,4] = 0.3218
$data[58,5] = 6448985.5
$data[59,0] = 45539.41666666666
$data[59,1] = 0.3218
$data[59,2] = 0.322
$data[59,3] = 0.315
$data[59,4] = 0.3175
$data[59,5] = 10307194.4
$data[60,0] = 45539.5
$data[60,1] = 0.3176
$data[60,2] = 0.3204
$data[60,3] = 0.3163
$data[60,4] = 0.318
$data[60,5] = 4845093.8
$data[61,0] = 45539.58333333334
$data[61,1] = 0.318
$data[61,2] = 0.3262
$data[61,3] = 0.3154
$data[61,4] = 0.326
$data[61,5] = 12542757.4
$data[62,0] = 45539.66666666666
$data[62,1] = 0.3259
$data[62,2] = 0.3293
$data[62,3] = 0.3223
$data[62,4] = 0.3223
$data[62,5] = 10636847.5
$data[63,0] = 45539.75
$data[63,1] = 0.3223
$data[63,2] = 0.3233
$data[63,3] = 0.321
$data[63,4] = 0.3223
$data[63,5] = 3525177

$ws.Range("A1169:F1232").Value = $data
